$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'60.440.56"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +3.35%  '
$ws.Range('D3').Value = "'2.327.12"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'544.95"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.44%  '
$ws.Range('D6').Value = "'131.17"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -1.47%  '
$ws.Range('D9').Value = "'2.324.62"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.14%  '
$ws.Range('E10').Value = '  +0.50%  '
$ws.Range('E11').Value = '  +0.50%  '
$ws.Range('E12').Value = '  +0.22%  '
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').Value = "'23.61"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').Value = "'60.441.59"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.41%  '
$ws.Range('D16').Value = "'2.740.44"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.08%  '
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').Value = "'2.317.01"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.58%  '
$ws.Range('D19').Value = "'10.58"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('D21').Value = "'314.26"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.60%  '
$ws.Range('D22').Value = "'6.66"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.49%  '
$ws.Range('D23').Value = "'1.00"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').Value = "'63.89"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.37%  '
$ws.Range('E25').Value = '  +2.54%  '
$ws.Range('D26').Value = "'0.995"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.62%  '
$ws.Range('D27').Value = "'7.86"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.13%  '
$ws.Range('E28').Value = '  +4.79%  '
$ws.Range('E29').Value = '  +9.24%  '
$ws.Range('D30').Value = "'173.52"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.50%  '
$ws.Range('E31').Value = '  +1.94%  '
$ws.Range('E32').Value = '  +1.01%  '
$ws.Range('E33').Value = '  +1.96%  '
$ws.Range('D34').Value = "'1.38"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +10.88%  '
$ws.Range('E35').Value = '  -0.75%  '
$ws.Range('D37').Value = "'17.83"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('D38').Value = "'1.00"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('D39').Value = "'4.07"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.76%  '
$ws.Range('D40').Value = "'322.36"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +11.02%  '
$ws.Range('E41').Value = '  +2.06%  '
$ws.Range('D42').Value = "'37.97"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.99%  '
$ws.Range('D43').Value = "'137.87"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.24%  '
$ws.Range('E44').Value = '  +1.01%  '
$ws.Range('E45').Value = '  -1.18%  '
$ws.Range('D46').Value = "'19.10"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +4.62%  '
$ws.Range('D47').Value = "'0.563"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.30%  '
$ws.Range('D48').Value = "'0.0495"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.35%  '
$ws.Range('E49').Value = '  +1.22%  '
$ws.Range('D50').Value = "'0.0₆0214"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +16.46%  '
$ws.Range('E51').Value = '  +0.62%  '
